$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.458797812461853
$ws.Range("B1").Value = 1.995802879333496
$ws.Range("C1").Value = 5.921448230743408
$ws.Range("D1").Value = 1.661562919616699
$ws.Range("E1").Value = 0.8170109391212463
